$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 8559.5
$ws.Range("I62").Value = 17000.572
$ws.Range("J62").Value = 3187.9092
$ws.Range("K62").Value = 17000.572
$ws.Range("L62").Value = 3187.9092
$ws.Range("M62").Value = -16376.572
$ws.Range("N62").Value = -4435.9092

$ws.Range("H65").Value = 8559.5
$ws.Range("I65").Value = 17000.572
$ws.Range("J65").Value = 3187.9092
$ws.Range("K65").Value = 85002.86
$ws.Range("L65").Value = 15939.546
$ws.Range("M65").Value = -81882.86
$ws.Range("N65").Value = -22179.546

$ws.Range("H129").Value = 1121.5143
$ws.Range("I129").Value = 415.2
$ws.Range("J129").Value = 1239.2333
$ws.Range("K129").Value = 1245.6
$ws.Range("L129").Value = 3717.699900000001
$ws.Range("M129").Value = 3754.4
$ws.Range("N129").Value = -13717.6999

$ws.Range("H132").Value = 3316.4546
$ws.Range("I132").Value = 3720.2222
$ws.Range("J132").Value = 1499.5
$ws.Range("K132").Value = 11160.6666
$ws.Range("L132").Value = 4498.5
$ws.Range("M132").Value = -8630.6666
$ws.Range("N132").Value = -9558.5

$ws.Range("H133").Value = 57800
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 57800
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 57800
$ws.Range("N133").Value = -67920
$ws.Range("M133").ClearContents()

$ws.Range("H134").Value = 49500
$ws.Range("J134").Value = 49500
$ws.Range("L134").Value = 49500
$ws.Range("N134").Value = -59640

$ws.Range("H137").Value = 1012.7857
$ws.Range("I137").Value = 825.4545000000001
$ws.Range("J137").Value = 1699.6666
$ws.Range("K137").Value = 2476.3635
$ws.Range("L137").Value = 5098.9998
$ws.Range("M137").Value = 73.63649999999961
$ws.Range("N137").Value = -10198.9998

$ws.Range("H138").Value = 3917.125
$ws.Range("I138").Value = 2464.3809
$ws.Range("J138").Value = 5522.7896
$ws.Range("K138").Value = 7393.1427
$ws.Range("L138").Value = 16568.3688
$ws.Range("M138").Value = -2253.1427
$ws.Range("N138").Value = -26848.3688

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 24329.75
$ws.Range("J55").Value = 24329.75
$ws.Range("L55").Value = 24329.75
$ws.Range("N55").Value = -24959.75

$ws.Range("H74").Value = 2083.0588
$ws.Range("I74").Value = 2083.0588
$ws.Range("K74").Value = 2083.0588
$ws.Range("M74").Value = -1209.0588

$ws.Range("H76").Value = 22000
$ws.Range("J76").Value = 22000
$ws.Range("L76").Value = 22000
$ws.Range("N76").Value = -22676

$ws.Range("H77").Value = 2083.0588
$ws.Range("I77").Value = 2083.0588
$ws.Range("K77").Value = 10415.294
$ws.Range("M77").Value = -6047.293999999998

$ws.Range("H79").Value = 22000
$ws.Range("J79").Value = 22000
$ws.Range("L79").Value = 22000
$ws.Range("N79").Value = -24340

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 48681.047
$ws.Range("I107").Value = 754.1177
$ws.Range("J107").Value = 252370.5
$ws.Range("K107").Value = 754.1177
$ws.Range("L107").Value = 252370.5
$ws.Range("M107").Value = 1165.8823
$ws.Range("N107").Value = -256210.5

$ws.Range("H134").Value = 2101.121
$ws.Range("I134").Value = 1512.0358
$ws.Range("J134").Value = 5400
$ws.Range("K134").Value = 4536.107400000001
$ws.Range("L134").Value = 16200
$ws.Range("M134").Value = -2001.107400000001
$ws.Range("N134").Value = -21270

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1563.4255
$ws.Range("I31").Value = 1086.2333
$ws.Range("J31").Value = 2405.5293
$ws.Range("K31").Value = 1086.2333
$ws.Range("L31").Value = 2405.5293
$ws.Range("M31").Value = -791.2333000000001
$ws.Range("N31").Value = -2995.5293

$ws.Range("H34").Value = 1563.4255
$ws.Range("I34").Value = 1086.2333
$ws.Range("J34").Value = 2405.5293
$ws.Range("K34").Value = 1086.2333
$ws.Range("L34").Value = 2405.5293
$ws.Range("M34").Value = -884.2333000000001
$ws.Range("N34").Value = -2809.5293

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 10785.714
$ws.Range("I87").Value = 4000
$ws.Range("J87").Value = 12636.363
$ws.Range("K87").Value = 12000
$ws.Range("L87").Value = 37909.089
$ws.Range("M87").Value = -10752
$ws.Range("N87").Value = -40405.089

$ws.Range("H90").Value = 10785.714
$ws.Range("I90").Value = 4000
$ws.Range("J90").Value = 12636.363
$ws.Range("K90").Value = 36000
$ws.Range("L90").Value = 113727.267
$ws.Range("M90").Value = -29760
$ws.Range("N90").Value = -126207.267

$ws.Range("H134").Value = 3323.0256
$ws.Range("I134").Value = 2149.9
$ws.Range("J134").Value = 4557.8945
$ws.Range("K134").Value = 6449.700000000001
$ws.Range("L134").Value = 13673.6835
$ws.Range("M134").Value = -1379.700000000001
$ws.Range("N134").Value = -23813.6835

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 8221.895
$ws.Range("I132").Value = 19687.334
$ws.Range("J132").Value = 2930.1538
$ws.Range("K132").Value = 59062.00199999999
$ws.Range("L132").Value = 8790.4614
$ws.Range("M132").Value = -56532.00199999999
$ws.Range("N132").Value = -13850.4614

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4458.625
$ws.Range("I132").Value = 4388.885
$ws.Range("J132").Value = 4760.8335
$ws.Range("K132").Value = 13166.655
$ws.Range("L132").Value = 14282.5005
$ws.Range("M132").Value = -10636.655
$ws.Range("N132").Value = -19342.5005

$ws.Range("H134").Value = 50000
$ws.Range("J134").Value = 50000
$ws.Range("L134").Value = 50000
$ws.Range("N134").Value = -60140

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws.Range("H136").Value = 4962.1787
$ws.Range("I136").Value = 3681.3333
$ws.Range("J136").Value = 12647.25
$ws.Range("K136").Value = 11043.9999
$ws.Range("L136").Value = 37941.75
$ws.Range("M136").Value = -8493.999899999999
$ws.Range("N136").Value = -43041.75

$ws.Range("H137").Value = 35000
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws.Range("H138").Value = 69800
$ws.Range("J138").Value = 69800
$ws.Range("L138").Value = 69800
$ws.Range("N138").Value = -80080

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws.Range("H141").Value = 60000
$ws.Range("J141").Value = 60000
$ws.Range("L141").Value = 60000
$ws.Range("N141").Value = -70360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").ClearContents()

$ws.Range("H132").Value = 14708555
$ws.Range("I132").Value = 18519560
$ws.Range("J132").Value = 8966.643
$ws.Range("K132").Value = 55558680
$ws.Range("L132").Value = 26899.929
$ws.Range("M132").Value = -55556150
$ws.Range("N132").Value = -31959.929
